$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72, shifting existing rows 72-80 down to 73-81
$ws.Rows.Item(72).Insert()

# Populate row 72
$ws.Cells.Item(72, 1).NumberFormat = "@"
$ws.Cells.Item(72, 1).Value = '7702'
$ws.Cells.Item(72, 1).ClearFormats()
$ws.Cells.Item(72, 2).NumberFormat = "@"
$ws.Cells.Item(72, 2).Value = '10/13/2025'
$ws.Cells.Item(72, 2).ClearFormats()
$ws.Cells.Item(72, 3).NumberFormat = "@"
$ws.Cells.Item(72, 3).Value = 'PAZ, MARCOS 3601'
$ws.Cells.Item(72, 3).ClearFormats()
$ws.Cells.Item(72, 4).NumberFormat = "@"
$ws.Cells.Item(72, 4).Value = '11'
$ws.Cells.Item(72, 4).ClearFormats()
$ws.Cells.Item(72, 5).NumberFormat = "@"
$ws.Cells.Item(72, 5).Value = '810333025'
$ws.Cells.Item(72, 5).ClearFormats()
$ws.Cells.Item(72, 6).NumberFormat = "@"
$ws.Cells.Item(72, 6).Value = 'NEW'
$ws.Cells.Item(72, 6).ClearFormats()
$ws.Cells.Item(72, 7).NumberFormat = "@"
$ws.Cells.Item(72, 7).Value = 'Pendiente'
$ws.Cells.Item(72, 7).ClearFormats()
$ws.Cells.Item(72, 8).NumberFormat = "@"
$ws.Cells.Item(72, 8).Value = 'Picada'
$ws.Cells.Item(72, 8).ClearFormats()
$ws.Cells.Item(72, 9).Value = 1
$ws.Cells.Item(72, 10).NumberFormat = "@"
$ws.Cells.Item(72, 10).Value = 'Cambio'
$ws.Cells.Item(72, 10).ClearFormats()
$ws.Cells.Item(72, 11).NumberFormat = "@"
$ws.Cells.Item(72, 11).Value = 'Sin equipos'
$ws.Cells.Item(72, 11).ClearFormats()
$ws.Cells.Item(72, 12).NumberFormat = "@"
$ws.Cells.Item(72, 12).Value = 'Pasante'
$ws.Cells.Item(72, 12).ClearFormats()
$ws.Cells.Item(72, 13).Value = -58.515009
$ws.Cells.Item(72, 14).Value = -34.606519
$ws.Cells.Item(72, 15).NumberFormat = "@"
$ws.Cells.Item(72, 15).Value = 'Devoto'
$ws.Cells.Item(72, 15).ClearFormats()
$ws.Cells.Item(72, 16).NumberFormat = "@"
$ws.Cells.Item(72, 16).Value = 'Capital Norte'
$ws.Cells.Item(72, 16).ClearFormats()
$ws.Cells.Item(72, 17).NumberFormat = "@"
$ws.Cells.Item(72, 17).Value = 'DEV-F'
$ws.Cells.Item(72, 17).ClearFormats()
$ws.Cells.Item(72, 18).NumberFormat = "@"
$ws.Cells.Item(72, 18).Value = 'ARATO-25058.PO.1DEV'
$ws.Cells.Item(72, 18).ClearFormats()

# Append new rows 82-85 at the end of the data
# Populate row 82
$ws.Cells.Item(82, 1).NumberFormat = "@"
$ws.Cells.Item(82, 1).Value = '7697'
$ws.Cells.Item(82, 1).ClearFormats()
$ws.Cells.Item(82, 2).NumberFormat = "@"
$ws.Cells.Item(82, 2).Value = '10/30/2025'
$ws.Cells.Item(82, 2).ClearFormats()
$ws.Cells.Item(82, 3).NumberFormat = "@"
$ws.Cells.Item(82, 3).Value = 'CONDE 4334'
$ws.Cells.Item(82, 3).ClearFormats()
$ws.Cells.Item(82, 4).NumberFormat = "@"
$ws.Cells.Item(82, 4).Value = '12'
$ws.Cells.Item(82, 4).ClearFormats()
$ws.Cells.Item(82, 5).NumberFormat = "@"
$ws.Cells.Item(82, 5).Value = '810487016'
$ws.Cells.Item(82, 5).ClearFormats()
$ws.Cells.Item(82, 6).NumberFormat = "@"
$ws.Cells.Item(82, 6).Value = 'NEW'
$ws.Cells.Item(82, 6).ClearFormats()
$ws.Cells.Item(82, 7).NumberFormat = "@"
$ws.Cells.Item(82, 7).Value = 'Pendiente'
$ws.Cells.Item(82, 7).ClearFormats()
$ws.Cells.Item(82, 8).NumberFormat = "@"
$ws.Cells.Item(82, 8).Value = 'Poste para cambiar o desmontar ver con inspector'
$ws.Cells.Item(82, 8).ClearFormats()
$ws.Cells.Item(82, 9).Value = 1
$ws.Cells.Item(82, 10).NumberFormat = "@"
$ws.Cells.Item(82, 10).Value = 'Cambio'
$ws.Cells.Item(82, 10).ClearFormats()
$ws.Cells.Item(82, 11).NumberFormat = "@"
$ws.Cells.Item(82, 11).Value = 'Sin equipos'
$ws.Cells.Item(82, 11).ClearFormats()
$ws.Cells.Item(82, 12).NumberFormat = "@"
$ws.Cells.Item(82, 12).Value = 'Poste'
$ws.Cells.Item(82, 12).ClearFormats()
$ws.Cells.Item(82, 13).Value = -58.481509
$ws.Cells.Item(82, 14).Value = -34.547874
$ws.Cells.Item(82, 15).NumberFormat = "@"
$ws.Cells.Item(82, 15).Value = 'Saavedra'
$ws.Cells.Item(82, 15).ClearFormats()
$ws.Cells.Item(82, 16).NumberFormat = "@"
$ws.Cells.Item(82, 16).Value = 'Capital Norte'
$ws.Cells.Item(82, 16).ClearFormats()
$ws.Cells.Item(82, 17).NumberFormat = "@"
$ws.Cells.Item(82, 17).Value = 'COG-P'
$ws.Cells.Item(82, 17).ClearFormats()
$ws.Cells.Item(82, 18).NumberFormat = "@"
$ws.Cells.Item(82, 18).Value = 'Fuera de Poligono OVL'
$ws.Cells.Item(82, 18).ClearFormats()

# Populate row 83
$ws.Cells.Item(83, 1).NumberFormat = "@"
$ws.Cells.Item(83, 1).Value = '7718'
$ws.Cells.Item(83, 1).ClearFormats()
$ws.Cells.Item(83, 2).NumberFormat = "@"
$ws.Cells.Item(83, 2).Value = '10/30/2025'
$ws.Cells.Item(83, 2).ClearFormats()
$ws.Cells.Item(83, 3).NumberFormat = "@"
$ws.Cells.Item(83, 3).Value = 'DORREGO 2293'
$ws.Cells.Item(83, 3).ClearFormats()
$ws.Cells.Item(83, 4).NumberFormat = "@"
$ws.Cells.Item(83, 4).Value = '14'
$ws.Cells.Item(83, 4).ClearFormats()
$ws.Cells.Item(83, 5).NumberFormat = "@"
$ws.Cells.Item(83, 5).Value = '810487028'
$ws.Cells.Item(83, 5).ClearFormats()
$ws.Cells.Item(83, 6).NumberFormat = "@"
$ws.Cells.Item(83, 6).Value = 'NEW'
$ws.Cells.Item(83, 6).ClearFormats()
$ws.Cells.Item(83, 7).NumberFormat = "@"
$ws.Cells.Item(83, 7).Value = 'Pendiente'
$ws.Cells.Item(83, 7).ClearFormats()
$ws.Cells.Item(83, 8).NumberFormat = "@"
$ws.Cells.Item(83, 8).Value = 'Picada'
$ws.Cells.Item(83, 8).ClearFormats()
$ws.Cells.Item(83, 9).Value = 1
$ws.Cells.Item(83, 10).NumberFormat = "@"
$ws.Cells.Item(83, 10).Value = 'Cambio'
$ws.Cells.Item(83, 10).ClearFormats()
$ws.Cells.Item(83, 11).NumberFormat = "@"
$ws.Cells.Item(83, 11).Value = 'Sin equipos'
$ws.Cells.Item(83, 11).ClearFormats()
$ws.Cells.Item(83, 12).NumberFormat = "@"
$ws.Cells.Item(83, 12).Value = 'Pasante'
$ws.Cells.Item(83, 12).ClearFormats()
$ws.Cells.Item(83, 13).Value = -58.437895
$ws.Cells.Item(83, 14).Value = -34.57696
$ws.Cells.Item(83, 15).NumberFormat = "@"
$ws.Cells.Item(83, 15).Value = 'Palermo'
$ws.Cells.Item(83, 15).ClearFormats()
$ws.Cells.Item(83, 16).NumberFormat = "@"
$ws.Cells.Item(83, 16).Value = 'Capital Sur'
$ws.Cells.Item(83, 16).ClearFormats()
$ws.Cells.Item(83, 17).NumberFormat = "@"
$ws.Cells.Item(83, 17).Value = 'ATH-B'
$ws.Cells.Item(83, 17).ClearFormats()
$ws.Cells.Item(83, 18).NumberFormat = "@"
$ws.Cells.Item(83, 18).Value = 'Fuera de Poligono OVL'
$ws.Cells.Item(83, 18).ClearFormats()

# Populate row 84
$ws.Cells.Item(84, 1).NumberFormat = "@"
$ws.Cells.Item(84, 1).Value = '7725'
$ws.Cells.Item(84, 1).ClearFormats()
$ws.Cells.Item(84, 2).NumberFormat = "@"
$ws.Cells.Item(84, 2).Value = '10/30/2025'
$ws.Cells.Item(84, 2).ClearFormats()
$ws.Cells.Item(84, 3).NumberFormat = "@"
$ws.Cells.Item(84, 3).Value = 'CHARCAS 4028'
$ws.Cells.Item(84, 3).ClearFormats()
$ws.Cells.Item(84, 4).NumberFormat = "@"
$ws.Cells.Item(84, 4).Value = '14'
$ws.Cells.Item(84, 4).ClearFormats()
$ws.Cells.Item(84, 5).NumberFormat = "@"
$ws.Cells.Item(84, 5).Value = '810487030'
$ws.Cells.Item(84, 5).ClearFormats()
$ws.Cells.Item(84, 6).NumberFormat = "@"
$ws.Cells.Item(84, 6).Value = 'NEW'
$ws.Cells.Item(84, 6).ClearFormats()
$ws.Cells.Item(84, 7).NumberFormat = "@"
$ws.Cells.Item(84, 7).Value = 'Pendiente'
$ws.Cells.Item(84, 7).ClearFormats()
$ws.Cells.Item(84, 8).NumberFormat = "@"
$ws.Cells.Item(84, 8).Value = 'Picada'
$ws.Cells.Item(84, 8).ClearFormats()
$ws.Cells.Item(84, 9).Value = 1
$ws.Cells.Item(84, 10).NumberFormat = "@"
$ws.Cells.Item(84, 10).Value = 'Cambio'
$ws.Cells.Item(84, 10).ClearFormats()
$ws.Cells.Item(84, 11).NumberFormat = "@"
$ws.Cells.Item(84, 11).Value = 'Sin equipos'
$ws.Cells.Item(84, 11).ClearFormats()
$ws.Cells.Item(84, 12).NumberFormat = "@"
$ws.Cells.Item(84, 12).Value = 'Pasante'
$ws.Cells.Item(84, 12).ClearFormats()
$ws.Cells.Item(84, 13).Value = -58.420243
$ws.Cells.Item(84, 14).Value = -34.585909
$ws.Cells.Item(84, 15).NumberFormat = "@"
$ws.Cells.Item(84, 15).Value = 'Palermo'
$ws.Cells.Item(84, 15).ClearFormats()
$ws.Cells.Item(84, 16).NumberFormat = "@"
$ws.Cells.Item(84, 16).Value = 'Capital Sur'
$ws.Cells.Item(84, 16).ClearFormats()
$ws.Cells.Item(84, 17).NumberFormat = "@"
$ws.Cells.Item(84, 17).Value = 'VCR-O'
$ws.Cells.Item(84, 17).ClearFormats()
$ws.Cells.Item(84, 18).NumberFormat = "@"
$ws.Cells.Item(84, 18).Value = 'Fuera de Poligono OVL'
$ws.Cells.Item(84, 18).ClearFormats()

# Populate row 85
$ws.Cells.Item(85, 1).NumberFormat = "@"
$ws.Cells.Item(85, 1).Value = '7726'
$ws.Cells.Item(85, 1).ClearFormats()
$ws.Cells.Item(85, 2).NumberFormat = "@"
$ws.Cells.Item(85, 2).Value = '10/30/2025'
$ws.Cells.Item(85, 2).ClearFormats()
$ws.Cells.Item(85, 3).NumberFormat = "@"
$ws.Cells.Item(85, 3).Value = 'GUEMES 4205'
$ws.Cells.Item(85, 3).ClearFormats()
$ws.Cells.Item(85, 4).NumberFormat = "@"
$ws.Cells.Item(85, 4).Value = '14'
$ws.Cells.Item(85, 4).ClearFormats()
$ws.Cells.Item(85, 5).NumberFormat = "@"
$ws.Cells.Item(85, 5).Value = '810487032'
$ws.Cells.Item(85, 5).ClearFormats()
$ws.Cells.Item(85, 6).NumberFormat = "@"
$ws.Cells.Item(85, 6).Value = 'NEW'
$ws.Cells.Item(85, 6).ClearFormats()
$ws.Cells.Item(85, 7).NumberFormat = "@"
$ws.Cells.Item(85, 7).Value = 'Pendiente'
$ws.Cells.Item(85, 7).ClearFormats()
$ws.Cells.Item(85, 8).NumberFormat = "@"
$ws.Cells.Item(85, 8).Value = 'Picada'
$ws.Cells.Item(85, 8).ClearFormats()
$ws.Cells.Item(85, 9).Value = 1
$ws.Cells.Item(85, 10).NumberFormat = "@"
$ws.Cells.Item(85, 10).Value = 'Cambio'
$ws.Cells.Item(85, 10).ClearFormats()
$ws.Cells.Item(85, 11).NumberFormat = "@"
$ws.Cells.Item(85, 11).Value = 'Sin equipos'
$ws.Cells.Item(85, 11).ClearFormats()
$ws.Cells.Item(85, 12).NumberFormat = "@"
$ws.Cells.Item(85, 12).Value = 'Pasante'
$ws.Cells.Item(85, 12).ClearFormats()
$ws.Cells.Item(85, 13).Value = -58.419744
$ws.Cells.Item(85, 14).Value = -34.584338
$ws.Cells.Item(85, 15).NumberFormat = "@"
$ws.Cells.Item(85, 15).Value = 'Palermo'
$ws.Cells.Item(85, 15).ClearFormats()
$ws.Cells.Item(85, 16).NumberFormat = "@"
$ws.Cells.Item(85, 16).Value = 'Capital Sur'
$ws.Cells.Item(85, 16).ClearFormats()
$ws.Cells.Item(85, 17).NumberFormat = "@"
$ws.Cells.Item(85, 17).Value = 'VCR-L'
$ws.Cells.Item(85, 17).ClearFormats()
$ws.Cells.Item(85, 18).NumberFormat = "@"
$ws.Cells.Item(85, 18).Value = 'Fuera de Poligono OVL'
$ws.Cells.Item(85, 18).ClearFormats()
